# Apply the edit described in the commit: "Removed confidential data and updated tests."
#
# This removes the "IgGI1H4N4S1" analyte (previously row 9) entirely, which shifts
# every row below it up by one, and also removes the trailing "test" row
# (previously the last row), leaving the new last row present but empty
# (keeping its formatting).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Make sure Sheet1 is the active sheet (it was tabSelected in the original file).
$ws.Activate()

# Delete the row that held "IgGI1H4N4S1" (row 9). This shifts rows 10-16 up by one,
# so what used to be row 16 ("test") becomes row 15.
$ws.Rows.Item(9).Delete() | Out-Null

# The former "test" entry is now in row 15 (the last used row). Clear its
# contents but keep the existing cell formatting/style in place.
$ws.Range("A15").ClearContents() | Out-Null

# Update the used range/selection to match the new last cell, mirroring what
# Excel records when a user's selection ends on the new bottom cell.
$ws.Range("F16").Select() | Out-Null
